$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 4353.5
$ws.Range("I33").Value = 93.75
$ws.Range("J33").Value = 12873
$ws.Range("K33").Value = 93.75
$ws.Range("L33").Value = 12873
$ws.Range("M33").Value = 135.25
$ws.Range("N33").Value = -13331

$ws.Range("H51").Value = 2911.25
$ws.Range("I51").Value = 2461.75
$ws.Range("J51").Value = 3360.75
$ws.Range("K51").Value = 2461.75
$ws.Range("L51").Value = 3360.75
$ws.Range("M51").Value = -1977.75
$ws.Range("N51").Value = -4328.75

$ws.Range("H80").Value = 2169.5278
$ws.Range("I80").Value = 659
$ws.Range("J80").Value = 3521.0527
$ws.Range("K80").Value = 1977
$ws.Range("L80").Value = 10563.1581
$ws.Range("M80").Value = -979
$ws.Range("N80").Value = -12559.1581

$ws.Range("H83").Value = 2169.5278
$ws.Range("I83").Value = 659
$ws.Range("J83").Value = 3521.0527
$ws.Range("K83").Value = 5931
$ws.Range("L83").Value = 31689.4743
$ws.Range("M83").Value = -939
$ws.Range("N83").Value = -41673.4743

$ws.Range("H88").Value = 1226678.8
$ws.Range("I88").Value = 1105.5555
$ws.Range("J88").Value = 1962022.6
$ws.Range("K88").Value = 1105.5555
$ws.Range("L88").Value = 1962022.6
$ws.Range("M88").Value = -699.5554999999999
$ws.Range("N88").Value = -1962834.6

$ws.Range("H91").Value = 1226678.8
$ws.Range("I91").Value = 1105.5555
$ws.Range("J91").Value = 1962022.6
$ws.Range("K91").Value = 1105.5555
$ws.Range("L91").Value = 1962022.6
$ws.Range("M91").Value = 298.4445000000001
$ws.Range("N91").Value = -1964830.6

$ws.Range("H96").Value = 1128
$ws.Range("I96").Value = 1110.4
$ws.Range("J96").Value = 1150
$ws.Range("K96").Value = 3331.2
$ws.Range("L96").Value = 3450
$ws.Range("M96").Value = -1958.2
$ws.Range("N96").Value = -6196

$ws.Range("H125").Value = 1562
$ws.Range("I125").Value = 1183
$ws.Range("J125").Value = 1751.5
$ws.Range("K125").Value = 10647
$ws.Range("L125").Value = 15763.5
$ws.Range("M125").Value = -8187
$ws.Range("N125").Value = -20683.5

$ws.Range("H132").Value = 1565.5
$ws.Range("I132").Value = 1452.5435
$ws.Range("K132").Value = 4357.6305
$ws.Range("M132").Value = -1827.6305

$ws.Range("H135").Value = 3345.2058
$ws.Range("I135").Value = 1399.5714
$ws.Range("J135").Value = 12424.833
$ws.Range("K135").Value = 12596.1426
$ws.Range("L135").Value = 111823.497
$ws.Range("M135").Value = -10061.1426
$ws.Range("N135").Value = -116893.497

$ws.Range("H139").Value = 62334
$ws.Range("I139").Value = 30000
$ws.Range("J139").Value = 70417.5
$ws.Range("K139").Value = 30000
$ws.Range("L139").Value = 70417.5
$ws.Range("M139").Value = -24860
$ws.Range("N139").Value = -80697.5

$ws.Range("H140").Value = 69674.14
$ws.Range("J140").Value = 69674.14
$ws.Range("L140").Value = 69674.14
$ws.Range("N140").Value = -80034.14

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2254.5417
$ws.Range("I61").Value = 2011.6111
$ws.Range("J61").Value = 2983.3333
$ws.Range("K61").Value = 2011.6111
$ws.Range("L61").Value = 2983.3333
$ws.Range("M61").Value = -1799.6111
$ws.Range("N61").Value = -3407.3333

$ws.Range("H122").Value = 1799.1875
$ws.Range("I122").Value = 1912.0541
$ws.Range("J122").Value = 1419.5454
$ws.Range("K122").Value = 5736.1623
$ws.Range("L122").Value = 4258.6362
$ws.Range("M122").Value = -3286.1623
$ws.Range("N122").Value = -9158.636200000001

$ws.Range("H132").Value = 1401.8312
$ws.Range("I132").Value = 1029.2667
$ws.Range("K132").Value = 3087.800099999999
$ws.Range("M132").Value = -557.8000999999995

$ws.Range("H136").Value = 2254.5417
$ws.Range("I136").Value = 2011.6111
$ws.Range("J136").Value = 2983.3333
$ws.Range("K136").Value = 6034.8333
$ws.Range("L136").Value = 8949.999899999999
$ws.Range("M136").Value = -3484.8333
$ws.Range("N136").Value = -14049.9999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 5474.9375
$ws.Range("I134").Value = 1054.2307
$ws.Range("J134").Value = 24631.334
$ws.Range("K134").Value = 3162.6921
$ws.Range("L134").Value = 73894.00199999999
$ws.Range("M134").Value = -627.6921000000002
$ws.Range("N134").Value = -78964.00199999999

$ws.Range("H140").Value = 77835.8
$ws.Range("I140").Value = 39999
$ws.Range("J140").Value = 87295
$ws.Range("K140").Value = 39999
$ws.Range("L140").Value = 87295
$ws.Range("M140").Value = -34819
$ws.Range("N140").Value = -97655

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1922.1666
$ws.Range("I132").Value = 1502.2609
$ws.Range("J132").Value = 3301.8572
$ws.Range("K132").Value = 4506.7827
$ws.Range("L132").Value = 9905.571599999999
$ws.Range("M132").Value = -1976.7827
$ws.Range("N132").Value = -14965.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 2210
$ws.Range("J38").Value = 2980
$ws.Range("L38").Value = 8940
$ws.Range("N38").Value = -9634

$ws.Range("H113").Value = 804.85
$ws.Range("I113").Value = 469.8846
$ws.Range("J113").Value = 922.5405
$ws.Range("K113").Value = 1409.6538
$ws.Range("L113").Value = 2767.6215
$ws.Range("M113").Value = 760.3462
$ws.Range("N113").Value = -7107.621499999999

$ws.Range("H132").Value = 1374.4482
$ws.Range("I132").Value = 1442.7858
$ws.Range("J132").Value = 1310.6666
$ws.Range("K132").Value = 12985.0722
$ws.Range("L132").Value = 11795.9994
$ws.Range("M132").Value = -10455.0722
$ws.Range("N132").Value = -16855.9994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8278752
$ws.Range("I70").Value = 11486078
$ws.Range("J70").Value = 7227.579
$ws.Range("K70").Value = 11486078
$ws.Range("L70").Value = 7227.579
$ws.Range("M70").Value = -11485808
$ws.Range("N70").Value = -7767.579

$ws.Range("H73").Value = 8278752
$ws.Range("I73").Value = 11486078
$ws.Range("J73").Value = 7227.579
$ws.Range("K73").Value = 11486078
$ws.Range("L73").Value = 7227.579
$ws.Range("M73").Value = -11485142
$ws.Range("N73").Value = -9099.579

$ws.Range("H113").Value = 66668450
$ws.Range("I113").Value = 1810.1666
$ws.Range("J113").Value = 111112870
$ws.Range("K113").Value = 1810.1666
$ws.Range("L113").Value = 111112870
$ws.Range("M113").Value = 359.8334
$ws.Range("N113").Value = -111117210

$ws.Range("H140").Value = 58780
$ws.Range("J140").Value = 58780
$ws.Range("L140").Value = 58780
$ws.Range("N140").Value = -69140

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H32").Value = 742
$ws.Range("I32").Value = 742
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 742
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -425
$ws.Range("N32").ClearContents()

$ws.Range("H122").Value = 2484.8823
$ws.Range("I122").Value = 2306.1
$ws.Range("J122").Value = 2740.2856
$ws.Range("K122").Value = 6918.299999999999
$ws.Range("L122").Value = 8220.856800000001
$ws.Range("M122").Value = -4468.299999999999
$ws.Range("N122").Value = -13120.8568

$ws.Range("H132").Value = 4200.0303
$ws.Range("I132").Value = 5303.5
$ws.Range("J132").Value = 2502.3845
$ws.Range("K132").Value = 15910.5
$ws.Range("L132").Value = 7507.1535
$ws.Range("M132").Value = -13380.5
$ws.Range("N132").Value = -12567.1535

$ws.Range("H139").Value = 79300
$ws.Range("J139").Value = 79300
$ws.Range("L139").Value = 79300
$ws.Range("N139").Value = -89580

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 69633.336
$ws.Range("I96").Value = 1114.2858
$ws.Range("J96").Value = 129587.5
$ws.Range("K96").Value = 1114.2858
$ws.Range("L96").Value = 129587.5
$ws.Range("M96").Value = 258.7141999999999
$ws.Range("N96").Value = -132333.5

$ws.Range("H132").Value = 16396219
$ws.Range("I132").Value = 19608812
$ws.Range("J132").Value = 11989.7
$ws.Range("K132").Value = 58826436
$ws.Range("L132").Value = 35969.10000000001
$ws.Range("M132").Value = -58823906
$ws.Range("N132").Value = -41029.10000000001

$ws.Range("H136").Value = 862.8148
$ws.Range("I136").Value = 705.75
$ws.Range("J136").Value = 1091.2727
$ws.Range("K136").Value = 2117.25
$ws.Range("L136").Value = 3273.8181
$ws.Range("M136").Value = 432.75
$ws.Range("N136").Value = -8373.8181

